$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 79.20169066666666
$ws.Cells.Item(2, 8).Value = 237.605072
$ws.Cells.Item(2, 9).Value = 0.1882902756436698
$ws.Cells.Item(2, 10).Value = 0.1882902756436699
$ws.Cells.Item(2, 13).Value = 1.370876333333333
$ws.Cells.Item(2, 14).Value = 4.112629
$ws.Cells.Item(2, 15).Value = 0.01103063309339269
$ws.Cells.Item(2, 16).Value = 0.01103063309339269
$ws.Cells.Item(2, 17).Value = 108.5757232949209
$ws.Cells.Item(2, 18).Value = 977.1815096542881
$ws.Cells.Item(2, 19).Value = 0.002076960945679097
$ws.Cells.Item(2, 20).Value = 0.002076960945679097
$ws.Cells.Item(3, 7).Value = 79.20169066666666
$ws.Cells.Item(3, 8).Value = 237.605072
$ws.Cells.Item(3, 9).Value = 0.1882902756436698
$ws.Cells.Item(3, 10).Value = 0.1882902756436699
$ws.Cells.Item(3, 15).Value = 0.7476219244149905
$ws.Cells.Item(3, 16).Value = 0.7476219244149904
$ws.Cells.Item(3, 17).Value = 7358.924053336606
$ws.Cells.Item(3, 18).Value = 66230.31648002946
$ws.Cells.Item(3, 19).Value = 0.1407699382253495
$ws.Cells.Item(3, 20).Value = 0.1407699382253495
$ws.Cells.Item(4, 7).Value = 79.20169066666666
$ws.Cells.Item(4, 8).Value = 237.605072
$ws.Cells.Item(4, 9).Value = 0.1882902756436698
$ws.Cells.Item(4, 10).Value = 0.1882902756436699
$ws.Cells.Item(4, 13).Value = 29.718484
$ws.Cells.Item(4, 14).Value = 89.155452
$ws.Cells.Item(4, 15).Value = 0.2391271080585153
$ws.Cells.Item(4, 16).Value = 0.2391271080585153
$ws.Cells.Item(4, 17).Value = 2353.754176850282
$ws.Cells.Item(4, 18).Value = 21183.78759165254
$ws.Cells.Item(4, 19).Value = 0.04502530909021148
$ws.Cells.Item(4, 20).Value = 0.04502530909021148
$ws.Cells.Item(5, 7).Value = 79.20169066666666
$ws.Cells.Item(5, 8).Value = 237.605072
$ws.Cells.Item(5, 9).Value = 0.1882902756436698
$ws.Cells.Item(5, 10).Value = 0.1882902756436699
$ws.Cells.Item(5, 13).Value = 0.275941
$ws.Cells.Item(5, 14).Value = 0.827823
$ws.Cells.Item(5, 15).Value = 0.002220334433101459
$ws.Cells.Item(5, 16).Value = 0.002220334433101458
$ws.Cells.Item(5, 17).Value = 21.85499372425067
$ws.Cells.Item(5, 18).Value = 196.694943518256
$ws.Cells.Item(5, 19).Value = 0.0004180673824298051
$ws.Cells.Item(5, 20).Value = 0.0004180673824298051
$ws.Cells.Item(6, 9).Value = 0.3031383606299999
$ws.Cells.Item(6, 10).Value = 0.3031383606299999
$ws.Cells.Item(6, 13).Value = 1.370876333333333
$ws.Cells.Item(6, 14).Value = 4.112629
$ws.Cells.Item(6, 15).Value = 0.01103063309339269
$ws.Cells.Item(6, 16).Value = 0.01103063309339269
$ws.Cells.Item(6, 17).Value = 174.8017344566744
$ws.Cells.Item(6, 18).Value = 1573.21561011007
$ws.Cells.Item(6, 19).Value = 0.003343808032642086
$ws.Cells.Item(6, 20).Value = 0.003343808032642086
$ws.Cells.Item(7, 9).Value = 0.3031383606299999
$ws.Cells.Item(7, 10).Value = 0.3031383606299999
$ws.Cells.Item(7, 15).Value = 0.7476219244149905
$ws.Cells.Item(7, 16).Value = 0.7476219244149904
$ws.Cells.Item(7, 19).Value = 0.226632884538206
$ws.Cells.Item(7, 20).Value = 0.2266328845382059
$ws.Cells.Item(8, 9).Value = 0.3031383606299999
$ws.Cells.Item(8, 10).Value = 0.3031383606299999
$ws.Cells.Item(8, 13).Value = 29.718484
$ws.Cells.Item(8, 14).Value = 89.155452
$ws.Cells.Item(8, 15).Value = 0.2391271080585153
$ws.Cells.Item(8, 16).Value = 0.2391271080585153
$ws.Cells.Item(8, 17).Value = 3789.431929276573
$ws.Cells.Item(8, 18).Value = 34104.88736348916
$ws.Cells.Item(8, 19).Value = 0.07248859951905119
$ws.Cells.Item(8, 20).Value = 0.07248859951905116
$ws.Cells.Item(9, 9).Value = 0.3031383606299999
$ws.Cells.Item(9, 10).Value = 0.3031383606299999
$ws.Cells.Item(9, 13).Value = 0.275941
$ws.Cells.Item(9, 14).Value = 0.827823
$ws.Cells.Item(9, 15).Value = 0.002220334433101459
$ws.Cells.Item(9, 16).Value = 0.002220334433101458
$ws.Cells.Item(9, 17).Value = 35.18549721434333
$ws.Cells.Item(9, 18).Value = 316.66947492909
$ws.Cells.Item(9, 19).Value = 0.0006730685401007165
$ws.Cells.Item(9, 20).Value = 0.0006730685401007164
$ws.Cells.Item(10, 7).Value = 128.6091306666667
$ws.Cells.Item(10, 8).Value = 385.827392
$ws.Cells.Item(10, 9).Value = 0.3057491381773125
$ws.Cells.Item(10, 10).Value = 0.3057491381773124
$ws.Cells.Item(10, 13).Value = 1.370876333333333
$ws.Cells.Item(10, 14).Value = 4.112629
$ws.Cells.Item(10, 15).Value = 0.01103063309339269
$ws.Cells.Item(10, 16).Value = 0.01103063309339269
$ws.Cells.Item(10, 17).Value = 176.3072134815076
$ws.Cells.Item(10, 18).Value = 1586.764921333568
$ws.Cells.Item(10, 19).Value = 0.003372606561854959
$ws.Cells.Item(10, 20).Value = 0.003372606561854957
$ws.Cells.Item(11, 7).Value = 128.6091306666667
$ws.Cells.Item(11, 8).Value = 385.827392
$ws.Cells.Item(11, 9).Value = 0.3057491381773125
$ws.Cells.Item(11, 10).Value = 0.3057491381773124
$ws.Cells.Item(11, 15).Value = 0.7476219244149905
$ws.Cells.Item(11, 16).Value = 0.7476219244149904
$ws.Cells.Item(11, 17).Value = 11949.55331351232
$ws.Cells.Item(11, 18).Value = 107545.9798216108
$ws.Cells.Item(11, 19).Value = 0.2285847590723472
$ws.Cells.Item(11, 20).Value = 0.2285847590723471
$ws.Cells.Item(12, 7).Value = 128.6091306666667
$ws.Cells.Item(12, 8).Value = 385.827392
$ws.Cells.Item(12, 9).Value = 0.3057491381773125
$ws.Cells.Item(12, 10).Value = 0.3057491381773124
$ws.Cells.Item(12, 13).Value = 29.718484
$ws.Cells.Item(12, 14).Value = 89.155452
$ws.Cells.Item(12, 15).Value = 0.2391271080585153
$ws.Cells.Item(12, 16).Value = 0.2391271080585153
$ws.Cells.Item(12, 17).Value = 3822.068391971243
$ws.Cells.Item(12, 18).Value = 34398.61552774119
$ws.Cells.Item(12, 19).Value = 0.07311290720372414
$ws.Cells.Item(12, 20).Value = 0.07311290720372411
$ws.Cells.Item(13, 7).Value = 128.6091306666667
$ws.Cells.Item(13, 8).Value = 385.827392
$ws.Cells.Item(13, 9).Value = 0.3057491381773125
$ws.Cells.Item(13, 10).Value = 0.3057491381773124
$ws.Cells.Item(13, 13).Value = 0.275941
$ws.Cells.Item(13, 14).Value = 0.827823
$ws.Cells.Item(13, 15).Value = 0.002220334433101459
$ws.Cells.Item(13, 16).Value = 0.002220334433101458
$ws.Cells.Item(13, 17).Value = 35.48853212529067
$ws.Cells.Item(13, 18).Value = 319.396789127616
$ws.Cells.Item(13, 19).Value = 0.0006788653393861826
$ws.Cells.Item(13, 20).Value = 0.0006788653393861824
$ws.Cells.Item(14, 7).Value = 85.31435366666666
$ws.Cells.Item(14, 8).Value = 255.943061
$ws.Cells.Item(14, 9).Value = 0.2028222255490178
$ws.Cells.Item(14, 10).Value = 0.2028222255490178
$ws.Cells.Item(14, 13).Value = 1.370876333333333
$ws.Cells.Item(14, 14).Value = 4.112629
$ws.Cells.Item(14, 15).Value = 0.01103063309339269
$ws.Cells.Item(14, 16).Value = 0.01103063309339269
$ws.Cells.Item(14, 17).Value = 116.9554283352632
$ws.Cells.Item(14, 18).Value = 1052.598855017369
$ws.Cells.Item(14, 19).Value = 0.002237257553216553
$ws.Cells.Item(14, 20).Value = 0.002237257553216552
$ws.Cells.Item(15, 7).Value = 85.31435366666666
$ws.Cells.Item(15, 8).Value = 255.943061
$ws.Cells.Item(15, 9).Value = 0.2028222255490178
$ws.Cells.Item(15, 10).Value = 0.2028222255490178
$ws.Cells.Item(15, 15).Value = 0.7476219244149905
$ws.Cells.Item(15, 16).Value = 0.7476219244149904
$ws.Cells.Item(15, 17).Value = 7926.874338261172
$ws.Cells.Item(15, 18).Value = 71341.86904435055
$ws.Cells.Item(15, 19).Value = 0.1516343425790879
$ws.Cells.Item(15, 20).Value = 0.1516343425790879
$ws.Cells.Item(16, 7).Value = 85.31435366666666
$ws.Cells.Item(16, 8).Value = 255.943061
$ws.Cells.Item(16, 9).Value = 0.2028222255490178
$ws.Cells.Item(16, 10).Value = 0.2028222255490178
$ws.Cells.Item(16, 13).Value = 29.718484
$ws.Cells.Item(16, 14).Value = 89.155452
$ws.Cells.Item(16, 15).Value = 0.2391271080585153
$ws.Cells.Item(16, 16).Value = 0.2391271080585153
$ws.Cells.Item(16, 17).Value = 2535.413254413174
$ws.Cells.Item(16, 18).Value = 22818.71928971857
$ws.Cells.Item(16, 19).Value = 0.04850029224552854
$ws.Cells.Item(16, 20).Value = 0.04850029224552853
$ws.Cells.Item(17, 7).Value = 85.31435366666666
$ws.Cells.Item(17, 8).Value = 255.943061
$ws.Cells.Item(17, 9).Value = 0.2028222255490178
$ws.Cells.Item(17, 10).Value = 0.2028222255490178
$ws.Cells.Item(17, 13).Value = 0.275941
$ws.Cells.Item(17, 14).Value = 0.827823
$ws.Cells.Item(17, 15).Value = 0.002220334433101459
$ws.Cells.Item(17, 16).Value = 0.002220334433101458
$ws.Cells.Item(17, 17).Value = 23.54172806513366
$ws.Cells.Item(17, 18).Value = 211.875552586203
$ws.Cells.Item(17, 19).Value = 0.0004503331711847546
$ws.Cells.Item(17, 20).Value = 0.0004503331711847545
